$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 446, pushing the existing
# rows 446-506 down to 448-508 (dimension grows from R506 to R508).
$ws.Rows.Item(446).Insert()
$ws.Rows.Item(446).Insert()

# New row 446: Papa / Patagonia / "1a (cosecha)" from Región de Los Lagos
$ws.Cells.Item(446, 1).Value = 7
$ws.Cells.Item(446, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(446, 3).Value = "Ñuble"
$ws.Cells.Item(446, 4).Value = 44984
$ws.Cells.Item(446, 5).Value = 16
$ws.Cells.Item(446, 6).Value = 100114001
$ws.Cells.Item(446, 7).Value = "Papa"
$ws.Cells.Item(446, 8).Value = "Patagonia"
$ws.Cells.Item(446, 9).Value = "1a (cosecha)"
$ws.Cells.Item(446, 10).Value = 240
$ws.Cells.Item(446, 11).Value = 11500
$ws.Cells.Item(446, 12).Value = 12000
$ws.Cells.Item(446, 13).Value = 11750
$ws.Cells.Item(446, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(446, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(446, 16).Value = 470
$ws.Cells.Item(446, 17).Value = 25
$ws.Cells.Item(446, 18).Value = "Hortaliza"

# New row 447: Papa / Patagonia / "2a (cosecha)" from Región de Los Lagos
$ws.Cells.Item(447, 1).Value = 7
$ws.Cells.Item(447, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(447, 3).Value = "Ñuble"
$ws.Cells.Item(447, 4).Value = 44984
$ws.Cells.Item(447, 5).Value = 16
$ws.Cells.Item(447, 6).Value = 100114001
$ws.Cells.Item(447, 7).Value = "Papa"
$ws.Cells.Item(447, 8).Value = "Patagonia"
$ws.Cells.Item(447, 9).Value = "2a (cosecha)"
$ws.Cells.Item(447, 10).Value = 150
$ws.Cells.Item(447, 11).Value = 10000
$ws.Cells.Item(447, 12).Value = 10000
$ws.Cells.Item(447, 13).Value = 10000
$ws.Cells.Item(447, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(447, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(447, 16).Value = 400
$ws.Cells.Item(447, 17).Value = 25
$ws.Cells.Item(447, 18).Value = "Hortaliza"
